$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 32 (pushes existing rows 32-62 down to 33-63),
# carrying over formatting (incl. the date number format) from the row below.
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the new weekly price record.
$ws.Range("A32").Value = 9
$ws.Range("B32").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C32").Value = "Metropolitana"
$ws.Range("D32").Value = 44447
$ws.Range("E32").Value = 13
$ws.Range("F32").Value = 100112005
$ws.Range("G32").Value = "Puerro"
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 106
$ws.Range("K32").Value = 7000
$ws.Range("L32").Value = 8000
$ws.Range("M32").Value = 7500
$ws.Range("N32").Value = "`$/paquete 20 unidades"
$ws.Range("O32").Value = "Provincia de Chacabuco"
$ws.Range("P32").Value = 375
$ws.Range("Q32").Value = 20
$ws.Range("R32").Value = "Hortaliza"
